$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Approved/Rejected" status for row 2 from "Rejected" to "Approved"
$ws.Range("I2").Value = "Approved"

# Clear the "ReasonToReject" cell for row 2 (was "same testcases are repeating ")
$ws.Range("J2").ClearContents()

# Update the active selection to J2
$ws.Range("J2").Select()
